$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15  (diff hunk @@ -1358,22 +1358,22 @@)
$ws.Range("H15").Value = 4203
$ws.Range("I15").Value = 4203
$ws.Range("K15").Value = 12609
$ws.Range("M15").Value = -12440
# row 86  (diff hunk @@ -4867,25 +4867,25 @@)
$ws.Range("H86").Value = 1843.875
$ws.Range("I86").Value = 1749.4
$ws.Range("J86").Value = 2001.3334
$ws.Range("K86").Value = 1749.4
$ws.Range("L86").Value = 2001.3334
$ws.Range("M86").Value = -626.4000000000001
$ws.Range("N86").Value = -4247.3334
# row 89  (diff hunk @@ -5017,25 +5017,25 @@)
$ws.Range("H89").Value = 1843.875
$ws.Range("I89").Value = 1749.4
$ws.Range("J89").Value = 2001.3334
$ws.Range("K89").Value = 8747
$ws.Range("L89").Value = 10006.667
$ws.Range("M89").Value = -3131
$ws.Range("N89").Value = -21238.667
# row 106  (diff hunk @@ -5877,22 +5877,22 @@)
$ws.Range("H106").Value = 3385
$ws.Range("I106").Value = 2272.0908
$ws.Range("K106").Value = 2272.0908
$ws.Range("M106").Value = -1641.0908
# row 132  (diff hunk @@ -7154,25 +7154,25 @@)
$ws.Range("H132").Value = 963.3953
$ws.Range("I132").Value = 863.15
$ws.Range("J132").Value = 2300
$ws.Range("K132").Value = 2589.45
$ws.Range("L132").Value = 6900
$ws.Range("M132").Value = -59.44999999999982
$ws.Range("N132").Value = -11960
# row 138  (diff hunk @@ -7451,25 +7451,25 @@)
$ws.Range("H138").Value = 3730.75
$ws.Range("I138").Value = 4770.9165
$ws.Range("J138").Value = 2170.5
$ws.Range("K138").Value = 14312.7495
$ws.Range("L138").Value = 6511.5
$ws.Range("M138").Value = -9172.749500000002
$ws.Range("N138").Value = -16791.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 74  (diff hunk @@ -11236,22 +11236,22 @@)
$ws.Range("H74").Value = 3302.8823
$ws.Range("I74").Value = 3319.2856
$ws.Range("K74").Value = 3319.2856
$ws.Range("M74").Value = -2445.2856
# row 77  (diff hunk @@ -11383,22 +11383,22 @@)
$ws.Range("H77").Value = 3302.8823
$ws.Range("I77").Value = 3319.2856
$ws.Range("K77").Value = 16596.428
$ws.Range("M77").Value = -12228.428
# row 132  (diff hunk @@ -14033,22 +14033,22 @@)
$ws.Range("H132").Value = 1154.238
$ws.Range("I132").Value = 846.6111
$ws.Range("K132").Value = 2539.8333
$ws.Range("M132").Value = -9.833299999999781

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 99  (diff hunk @@ -26111,22 +26111,22 @@)
$ws.Range("H99").Value = 2325
$ws.Range("I99").Value = 2193.75
$ws.Range("K99").Value = 2193.75
$ws.Range("M99").Value = -695.75
# row 126  (diff hunk @@ -27416,22 +27416,22 @@)
$ws.Range("H126").Value = 2325
$ws.Range("I126").Value = 2193.75
$ws.Range("K126").Value = 6581.25
$ws.Range("M126").Value = -4111.25
# row 129  (diff hunk @@ -27563,22 +27563,22 @@)
$ws.Range("H129").Value = 30000
$ws.Range("J129").Value = 30000
$ws.Range("L129").Value = 30000
$ws.Range("N129").Value = -40000
# row 134  (diff hunk @@ -27802,25 +27802,25 @@)
$ws.Range("H134").Value = 1947.44
$ws.Range("I134").Value = 1718.4286
$ws.Range("J134").Value = 3149.75
$ws.Range("K134").Value = 5155.2858
$ws.Range("L134").Value = 9449.25
$ws.Range("M134").Value = -2620.2858
$ws.Range("N134").Value = -14519.25

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 5  (diff hunk @@ -28444,22 +28444,22 @@)
$ws.Range("H5").Value = 499.33334
$ws.Range("I5").Value = 499.33334
$ws.Range("K5").Value = 1498.00002
$ws.Range("M5").Value = -1386.00002
# row 68  (diff hunk @@ -31558,22 +31558,25 @@)
$ws.Range("H68").Value = 1998
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 1997.6666
$ws.Range("K68").Value = 5997
$ws.Range("L68").Value = 5992.9998
$ws.Range("M68").Value = -5186
$ws.Range("N68").Value = -7614.9998
# row 71  (diff hunk @@ -31708,22 +31711,25 @@)
$ws.Range("H71").Value = 1998
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 1997.6666
$ws.Range("K71").Value = 17991
$ws.Range("L71").Value = 17978.9994
$ws.Range("M71").Value = -13935
$ws.Range("N71").Value = -26090.9994
# row 121  (diff hunk @@ -34191,25 +34197,25 @@)
$ws.Range("H121").Value = 496.66666
$ws.Range("J121").Value = 490
$ws.Range("L121").Value = 1470
$ws.Range("N121").Value = -4090
# row 122  (diff hunk @@ -34243,25 +34249,25 @@)
$ws.Range("H122").Value = 1054.8572
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 976.8
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 8791.199999999999
$ws.Range("M122").Value = -8800
$ws.Range("N122").Value = -13691.2
# row 131  (diff hunk @@ -34696,25 +34702,25 @@)
$ws.Range("H131").Value = 6182436.5
$ws.Range("J131").Value = 10318.187
$ws.Range("L131").Value = 30954.561
$ws.Range("N131").Value = -41034.561
# row 132  (diff hunk @@ -34748,25 +34754,25 @@)
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 1600
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 14400
$ws.Range("L132").Value = 28800
$ws.Range("M132").Value = -11870
$ws.Range("N132").Value = -33860
# row 135  (diff hunk @@ -34904,22 +34910,22 @@)
$ws.Range("H135").Value = 499.33334
$ws.Range("I135").Value = 499.33334
$ws.Range("K135").Value = 4494.00006
$ws.Range("M135").Value = -1959.00006

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 14  (diff hunk @@ -35950,22 +35956,22 @@)
$ws.Range("H14").Value = 2330000
$ws.Range("I14").Value = 2330000
$ws.Range("K14").Value = 2330000
$ws.Range("M14").Value = -2329832
# row 102  (diff hunk @@ -40151,25 +40157,25 @@)
$ws.Range("H102").Value = 2492.7646
$ws.Range("I102").Value = 3086.1428
$ws.Range("J102").Value = 2077.4
$ws.Range("K102").Value = 3086.1428
$ws.Range("L102").Value = 2077.4
$ws.Range("M102").Value = -1464.1428
$ws.Range("N102").Value = -5321.4
# row 122  (diff hunk @@ -41110,22 +41116,22 @@)
$ws.Range("H122").Value = 1885.5834
$ws.Range("I122").Value = 1721
$ws.Range("K122").Value = 5163
$ws.Range("M122").Value = -2713
# row 132  (diff hunk @@ -41594,25 +41600,25 @@)
$ws.Range("H132").Value = 3472.423
$ws.Range("I132").Value = 2566.4119
$ws.Range("J132").Value = 5183.778
$ws.Range("K132").Value = 7699.2357
$ws.Range("L132").Value = 15551.334
$ws.Range("M132").Value = -5169.2357
$ws.Range("N132").Value = -20611.334
# row 134  (diff hunk @@ -41692,22 +41698,22 @@)
$ws.Range("H134").Value = 44999.285
$ws.Range("J134").Value = 44999.285
$ws.Range("L134").Value = 134997.855
$ws.Range("N134").Value = -140067.855
# row 135  (diff hunk @@ -41741,22 +41747,22 @@)
$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93  (diff hunk @@ -46631,22 +46637,22 @@)
$ws.Range("H93").Value = 1111.7142
$ws.Range("I93").Value = 965
$ws.Range("K93").Value = 965
$ws.Range("M93").Value = 283
# row 132  (diff hunk @@ -48482,22 +48488,22 @@)
$ws.Range("H132").Value = 1797.375
$ws.Range("I132").Value = 1485.8
$ws.Range("K132").Value = 4457.4
$ws.Range("M132").Value = -1927.4
# row 136  (diff hunk @@ -48681,25 +48687,25 @@)
$ws.Range("H136").Value = 3440.2666
$ws.Range("I136").Value = 2357.6667
$ws.Range("J136").Value = 5966.3335
$ws.Range("K136").Value = 7073.000100000001
$ws.Range("L136").Value = 17899.0005
$ws.Range("M136").Value = -4523.000100000001
$ws.Range("N136").Value = -22999.0005

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 132  (diff hunk @@ -55337,25 +55343,25 @@)
$ws.Range("H132").Value = 2165.7
$ws.Range("I132").Value = 1719.5
$ws.Range("J132").Value = 2463.1667
$ws.Range("K132").Value = 5158.5
$ws.Range("L132").Value = 7389.500100000001
$ws.Range("M132").Value = -2628.5
$ws.Range("N132").Value = -12449.5001
# row 136  (diff hunk @@ -55530,25 +55536,25 @@)
$ws.Range("H136").Value = 2512.718
$ws.Range("I136").Value = 2046.0667
$ws.Range("J136").Value = 4068.2222
$ws.Range("K136").Value = 6138.2001
$ws.Range("L136").Value = 12204.6666
$ws.Range("M136").Value = -3588.2001
$ws.Range("N136").Value = -17304.6666
# row 137  (diff hunk @@ -55582,22 +55588,19 @@)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""
